$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '90.972.79'
$ws.Range('E2').Value = '  +1.62%  '

# Row 3
$ws.Range('D3').Value = '3.168.46'
$ws.Range('E3').Value = '  +3.26%  '

# Row 4
$ws.Range('E4').Value = '  -0.23%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.09'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +2.41%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '628.15'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +3.01%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.19'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +31.02%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.368'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +1.03%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.00'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.10%  '

# Row 10
$ws.Range('D10').Value = '3.164.55'
$ws.Range('E10').Value = '  +3.07%  '

# Row 11
$ws.Range('E11').Value = '  +14.38%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.203'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +7.79%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000248'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +4.19%  '

# Row 14
$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.43'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +10.01%  '

# Row 15
$ws.Range('B15').Value = 'Toncoin'
$ws.Range('C15').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.66'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +5.14%  '

# Row 16
$ws.Range('D16').Value = '90.811.42'
$ws.Range('E16').Value = '  +1.60%  '

# Row 17
$ws.Range('D17').Value = '3.751.40'
$ws.Range('E17').Value = '  +2.92%  '

# Row 18
$ws.Range('D18').Value = '3.168.84'
$ws.Range('E18').Value = '  +2.57%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.70'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +9.31%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0000216'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +3.84%  '

# Row 21
$ws.Range('E21').Value = '  +6.23%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '455.62'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +6.43%  '

# Row 23
$ws.Range('E23').Value = '  +10.37%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.21'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +4.61%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.97'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +9.87%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '91.96'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +8.72%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.18'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +3.40%  '

# Row 28
$ws.Range('D28').Value = '3.336.23'
$ws.Range('E28').Value = '  +3.08%  '

# Row 29
$ws.Range('E29').Value = '  -0.01%  '

# Row 30
$ws.Range('B30').Value = 'Cronos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.163'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.83%  '

# Row 31
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.29'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +12.67%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.00'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -8.93%  '

# Row 33
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '25.99'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +14.24%  '

# Row 34
$ws.Range('B34').Value = 'Bittensor'
$ws.Range('C34').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '525.94'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +4.25%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.189'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +31.93%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.68'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +2.60%  '

# Row 37
$ws.Range('E37').Value = '  +9.73%  '

# Row 38
$ws.Range('E38').Value = '  +7.48%  '

# Row 39
$ws.Range('E39').Value = '  +4.95%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.30'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +5.37%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0922'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +33.44%  '

# Row 42
$ws.Range('B42').Value = 'WhiteBITCoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '22.23'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.10%  '

# Row 43
$ws.Range('B43').Value = 'PolygonEcosystemToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.424'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +15.18%  '

# Row 44
$ws.Range('E44').Value = '  -0.43%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.95'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +6.87%  '

# Row 46
$ws.Range('E46').Value = '  -0.04%  '

# Row 47
$ws.Range('B47').Value = 'ImmutableX'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.36'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +12.20%  '

# Row 48
$ws.Range('B48').Value = 'Monero'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '146.87'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.65%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.52'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +10.37%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '44.75'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +3.03%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.657'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +12.16%  '
